$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(27,1).Value = "1k1r4/pp4p1/2p4p/2PnR2P/5B2/6P1/PP3r2/1K1R4 b"
$ws.Cells.Item(27,2).Value = "d5c3 b2c3 d8d1"
$ws.Cells.Item(27,3).Value = "Discovered Attack"
$ws.Cells.Item(27,4).Value = 893

$ws.Cells.Item(28,1).Value = "r5k1/ppp3pp/5p2/8/3br3/7P/1PP1B1P1/2K1R2R w"
$ws.Cells.Item(28,2).Value = "e2c4 g8f8 e1e4"
$ws.Cells.Item(28,3).Value = "Discovered Attack"
$ws.Cells.Item(28,4).Value = 643

$ws.Cells.Item(29,1).Value = "r3r1k1/p4ppn/2p4p/8/3q2b1/2NBR3/PPPQ2PP/5RK1 w"
$ws.Cells.Item(29,2).Value = "d3h7 g8h7 d2d4"
$ws.Cells.Item(29,3).Value = "Discovered Attack"
$ws.Cells.Item(29,4).Value = 822

$ws.Cells.Item(30,2).Value = "e5f3 g2f3 e8e1 d1e1 g4f3"
$ws.Cells.Item(30,1).Value = "r2qr1k1/pp3ppp/3p1b2/2pPn3/2P2Bb1/3P1NP1/PP3PBP/1R1QR1K1 b"
$ws.Cells.Item(30,3).Value = "Discovered Attack"
$ws.Cells.Item(30,4).Value = 1612

$ws.Cells.Item(31,1).Value = "5rbR/2p1p3/8/pBkPp1p1/P3P3/6PK/1P5P/8 b"
$ws.Cells.Item(31,2).Value = "g8e6 d5e6 f8h8"
$ws.Cells.Item(31,3).Value = "Discovered Attack"
$ws.Cells.Item(31,4).Value = 1006

$ws.Cells.Item(32,1).Value = "2kr4/1ppn1prp/nq2p3/p6Q/P7/1PN5/2PP2BP/2KR2R1 w"
$ws.Cells.Item(32,2).Value = "g2b7 b6b7 g1g7"
$ws.Cells.Item(32,3).Value = "Discovered Attack"
$ws.Cells.Item(32,4).Value = 1273

$ws.Cells.Item(33,1).Value = "2kb2r1/1p1b2q1/4p3/3pN3/pNpP1R1p/P1P3r1/1P2Q1PK/6R1 b"
$ws.Cells.Item(33,2).Value = "g3h3 g2h3 g7g1"
$ws.Cells.Item(33,3).Value = "Discovered Attack"
$ws.Cells.Item(33,4).Value = 1581

$ws.Cells.Item(34,1).Value = "8/p6Q/kp6/3pP3/PP1P4/2P3q1/3BKp2/5r2 w"
$ws.Cells.Item(34,2).Value = "b4b5 a6a5 h7a7"
$ws.Cells.Item(34,3).Value = "Deflection"
$ws.Cells.Item(34,4).Value = 1710

$ws.Cells.Item(35,2).Value = "d1e2 f2g3 e2e3"
$ws.Cells.Item(35,1).Value = "6k1/5pp1/p1p1b2p/8/4P3/P1QpBP2/1P3KPP/3q4 b"
$ws.Cells.Item(35,3).Value = "Deflection"
$ws.Cells.Item(35,4).Value = 970

$ws.Cells.Item(36,1).Value = "8/8/4k3/5RP1/r4pK1/7P/6P1/8 b"
$ws.Cells.Item(36,2).Value = "f4f3 g4f3 e6f5"
$ws.Cells.Item(36,3).Value = "Deflection"
$ws.Cells.Item(36,4).Value = 914


# Re-sort the whole table (A2:D36) ascending by Rating (column D), matching
# the table's existing sort order/convention.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$keyRange = $ws.Range("D2:D36")
$sortObj.SortFields.Add($keyRange) | Out-Null
$sortObj.SetRange($ws.Range("A1:D36"))
$sortObj.Header = 1
$sortObj.Apply()

# Widen columns A and C to fit the new (longer) board/theme text.
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null

# Restore the active-cell selection as it ended up after the edit.
$ws.Range("F12").Select() | Out-Null
